# Remove the redundant explicit "left" paragraph justification (w:jc val="left")
# from the numbered list item paragraphs. Word omits w:jc entirely when the
# paragraph alignment equals the (left) default, so re-asserting
# wdAlignParagraphLeft (0) on these paragraphs clears the now-unneeded
# <w:jc w:val="left"/> element while leaving numbering/indent untouched.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.ListFormat.ListType -ne 0) {
        $p.Range.ParagraphFormat.Alignment = 0
    }
}
